# Update the deadline for SA5 from October 06, 2024 to October 07, 2024
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B17").Value = "October 07, 2024"

# Update the selected cell to match the saved selection in the workbook
$ws.Range("B17").Select()
